# Weekly update: insert this week's two new "Brócoli" price rows (Primera /
# Segunda) at the top of the data block (row 676), pushing the rest of the
# historical rows down by two. This mirrors how the upstream consolidated
# sheet gets a new week's record prepended above the prior weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the first data row of the block (row 676),
# shifting rows 676:770 down to 678:772 and extending the used range to R772.
$ws.Range("A676:A677").EntireRow.Insert()

# Row 676 -- new "Primera" record for 2022-08-03 (serial 44776)
$ws.Range("A676").Value = 8
$ws.Range("B676").Value = "Terminal La Palmera de La Serena"
$ws.Range("C676").Value = "Coquimbo"
$ws.Range("D676").Value = 44776
$ws.Range("E676").Value = 4
$ws.Range("F676").Value = 100112023
$ws.Range("G676").Value = "Brócoli"
$ws.Range("H676").Value = "Sin especificar"
$ws.Range("I676").Value = "Primera"
$ws.Range("J676").Value = 2640
$ws.Range("K676").Value = 750
$ws.Range("L676").Value = 800
$ws.Range("M676").Value = 775
$ws.Range("N676").Value = "$/unidad"
$ws.Range("O676").Value = "Provincia del Elquí"
$ws.Range("P676").Value = 775
$ws.Range("Q676").Value = 1
$ws.Range("R676").Value = "Hortaliza"

# Row 677 -- new "Segunda" record for 2022-08-03 (serial 44776)
$ws.Range("A677").Value = 8
$ws.Range("B677").Value = "Terminal La Palmera de La Serena"
$ws.Range("C677").Value = "Coquimbo"
$ws.Range("D677").Value = 44776
$ws.Range("E677").Value = 4
$ws.Range("F677").Value = 100112023
$ws.Range("G677").Value = "Brócoli"
$ws.Range("H677").Value = "Sin especificar"
$ws.Range("I677").Value = "Segunda"
$ws.Range("J677").Value = 1520
$ws.Range("K677").Value = 650
$ws.Range("L677").Value = 700
$ws.Range("M677").Value = 675
$ws.Range("N677").Value = "$/unidad"
$ws.Range("O677").Value = "Provincia del Elquí"
$ws.Range("P677").Value = 675
$ws.Range("Q677").Value = 1
$ws.Range("R677").Value = "Hortaliza"
